$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 113 (hunk 0)
$ws.Range("H113").Value = 2605.0527
$ws.Range("I113").Value = 2435.4167
$ws.Range("J113").Value = 2895.8572
$ws.Range("K113").Value = 2435.4167
$ws.Range("L113").Value = 2895.8572
$ws.Range("M113").Value = 818.5832999999998
$ws.Range("N113").Value = -9403.8572

# Row 116 (hunk 1)
$ws.Range("H116").Value = 3312.0588
$ws.Range("I116").Value = 2858.75
$ws.Range("J116").Value = 4400
$ws.Range("K116").Value = 2858.75
$ws.Range("L116").Value = 4400
$ws.Range("M116").Value = 583.25
$ws.Range("N116").Value = -11284

# Row 132 (hunk 2)
$ws.Range("H132").Value = 4465784.5
$ws.Range("I132").Value = 5436399
$ws.Range("J132").Value = 958
$ws.Range("K132").Value = 16309197
$ws.Range("L132").Value = 2874
$ws.Range("M132").Value = -16306667
$ws.Range("N132").Value = -7934

# Row 135 (hunk 3)
$ws.Range("H135").Value = 1627.1515
$ws.Range("I135").Value = 682.4167
$ws.Range("J135").Value = 4146.4443
$ws.Range("K135").Value = 6141.7503
$ws.Range("L135").Value = 37317.9987
$ws.Range("M135").Value = -3606.7503
$ws.Range("N135").Value = -42387.9987

# Row 137 (hunk 4)
$ws.Range("H137").Value = 3008.827
$ws.Range("I137").Value = 3443.3794
$ws.Range("J137").Value = 2460.913
$ws.Range("K137").Value = 10330.1382
$ws.Range("L137").Value = 7382.739
$ws.Range("M137").Value = -7780.138199999999
$ws.Range("N137").Value = -12482.739

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (hunk 5)
$ws.Range("H61").Value = 765.12067
$ws.Range("I61").Value = 771.0179000000001
$ws.Range("K61").Value = 771.0179000000001
$ws.Range("M61").Value = -559.0179000000001

# Row 132 (hunk 6)
$ws.Range("H132").Value = 1679.6305
$ws.Range("I132").Value = 1622.093
$ws.Range("J132").Value = 2504.3333
$ws.Range("K132").Value = 4866.279
$ws.Range("L132").Value = 7512.999899999999
$ws.Range("M132").Value = -2336.279
$ws.Range("N132").Value = -12572.9999

# Row 136 (hunk 7)
$ws.Range("H136").Value = 765.12067
$ws.Range("I136").Value = 771.0179000000001
$ws.Range("K136").Value = 2313.0537
$ws.Range("M136").Value = 236.9462999999996

$ws = $wb.Worksheets.Item("BSM")
# Row 22 (hunk 8)
$ws.Range("H22").Value = 365.55554
$ws.Range("I22").Value = 380
$ws.Range("J22").Value = 354
$ws.Range("K22").Value = 380
$ws.Range("L22").Value = 354
$ws.Range("M22").Value = -207
$ws.Range("N22").Value = -700

# Row 55 (hunk 9)
$ws.Range("H55").Value = 44926.668
$ws.Range("J55").Value = 44926.668
$ws.Range("L55").Value = 44926.668
$ws.Range("N55").Value = -45472.668

# Row 99 (hunk 10)
$ws.Range("H99").Value = 1037.762
$ws.Range("I99").Value = 978.36365
$ws.Range("J99").Value = 1255.5555
$ws.Range("K99").Value = 978.36365
$ws.Range("L99").Value = 1255.5555
$ws.Range("M99").Value = 519.63635
$ws.Range("N99").Value = -4251.5555

# Row 134 (hunk 11)
$ws.Range("H134").Value = 1214.25
$ws.Range("I134").Value = 1073.4
$ws.Range("J134").Value = 1449
$ws.Range("K134").Value = 3220.2
$ws.Range("L134").Value = 4347
$ws.Range("M134").Value = -685.2000000000003
$ws.Range("N134").Value = -9417

$ws = $wb.Worksheets.Item("CRP")
# Row 99 (hunk 12)
$ws.Range("H99").Value = 1750
$ws.Range("I99").Value = 1500
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -2
$ws.Range("N99").Value = -5496

# Row 122 (hunk 13)
$ws.Range("H122").Value = 2069.2
$ws.Range("I122").Value = 2069.2
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6207.599999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3757.599999999999
$ws.Range("N122").ClearContents()

# Row 123 (hunk 14)
$ws.Range("H123").Value = 30650
$ws.Range("J123").Value = 30650
$ws.Range("L123").Value = 30650
$ws.Range("N123").Value = -40450

# Row 126 (hunk 15)
$ws.Range("H126").Value = 1750
$ws.Range("I126").Value = 1500
$ws.Range("J126").Value = 2500
$ws.Range("K126").Value = 4500
$ws.Range("L126").Value = 7500
$ws.Range("M126").Value = -2030
$ws.Range("N126").Value = -12440

# Row 127 (hunk 16)
$ws.Range("H127").Value = 30740
$ws.Range("J127").Value = 30740
$ws.Range("L127").Value = 30740
$ws.Range("N127").Value = -40660

# Row 132 (hunk 17)
$ws.Range("H132").Value = 1043.6444
$ws.Range("I132").Value = 772.2941
$ws.Range("J132").Value = 1882.3636
$ws.Range("K132").Value = 2316.8823
$ws.Range("L132").Value = 5647.0908
$ws.Range("M132").Value = 213.1177000000002
$ws.Range("N132").Value = -10707.0908

# Row 134 (hunk 18)
$ws.Range("H134").Value = 1371.0435
$ws.Range("I134").Value = 1002.2857
$ws.Range("K134").Value = 3006.8571
$ws.Range("M134").Value = -471.8571000000002

$ws = $wb.Worksheets.Item("CUL")
# Row 131 (hunk 19)
$ws.Range("H131").Value = 837.6042
$ws.Range("I131").Value = 399.4
$ws.Range("J131").Value = 888.55817
$ws.Range("K131").Value = 1198.2
$ws.Range("L131").Value = 2665.67451
$ws.Range("M131").Value = 3841.8
$ws.Range("N131").Value = -12745.67451

$ws = $wb.Worksheets.Item("GSM")
# Row 32 (hunk 20)
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").ClearContents()

# Row 42 (hunk 21)
$ws.Range("H42").Value = 14745.454
$ws.Range("J42").Value = 14745.454
$ws.Range("L42").Value = 14745.454
$ws.Range("N42").Value = -15715.454

# Row 45 (hunk 22)
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = 0

# Row 51 (hunk 23)
$ws.Range("H51").Value = 39040
$ws.Range("I51").Value = 30000
$ws.Range("J51").Value = 41300
$ws.Range("K51").Value = 30000
$ws.Range("L51").Value = 41300
$ws.Range("M51").Value = -29491
$ws.Range("N51").Value = -42318

# Row 115 (hunk 24)
$ws.Range("H115").Value = 14745.454
$ws.Range("J115").Value = 14745.454
$ws.Range("L115").Value = 14745.454
$ws.Range("N115").Value = -17095.454

# Row 122 (hunk 25)
$ws.Range("H122").Value = 2248
$ws.Range("I122").Value = 2248
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 6744
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4294
$ws.Range("N122").ClearContents()

# Row 132 (hunk 26)
$ws.Range("H132").Value = 2782.348
$ws.Range("I132").Value = 2517.625
$ws.Range("J132").Value = 3387.4285
$ws.Range("K132").Value = 7552.875
$ws.Range("L132").Value = 10162.2855
$ws.Range("M132").Value = -5022.875
$ws.Range("N132").Value = -15222.2855

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (hunk 27)
$ws.Range("H7").Value = 2334.889
$ws.Range("I7").Value = 1051
$ws.Range("J7").Value = 3362
$ws.Range("K7").Value = 1051
$ws.Range("L7").Value = 3362
$ws.Range("M7").Value = -939
$ws.Range("N7").Value = -3586

# Row 93 (hunk 28)
$ws.Range("H93").Value = 875.9286
$ws.Range("I93").Value = 881.7692
$ws.Range("J93").Value = 800
$ws.Range("K93").Value = 881.7692
$ws.Range("L93").Value = 800
$ws.Range("M93").Value = 366.2308
$ws.Range("N93").Value = -3296

# Row 126 (hunk 29)
$ws.Range("H126").Value = 2334.889
$ws.Range("I126").Value = 1051
$ws.Range("J126").Value = 3362
$ws.Range("K126").Value = 3153
$ws.Range("L126").Value = 10086
$ws.Range("M126").Value = -683
$ws.Range("N126").Value = -15026

$ws = $wb.Worksheets.Item("WVR")
# Row 75 (hunk 30)
$ws.Range("H75").Value = 28000
$ws.Range("J75").Value = 28000
$ws.Range("L75").Value = 28000
$ws.Range("N75").Value = -29872

# Row 78 (hunk 31)
$ws.Range("H78").Value = 28000
$ws.Range("J78").Value = 28000
$ws.Range("L78").Value = 84000
$ws.Range("N78").Value = -93360

# Row 132 (hunk 32)
$ws.Range("H132").Value = 1333.5897
$ws.Range("I132").Value = 873.2632
$ws.Range("J132").Value = 1770.9
$ws.Range("K132").Value = 2619.7896
$ws.Range("L132").Value = 5312.700000000001
$ws.Range("M132").Value = -89.78960000000006
$ws.Range("N132").Value = -10372.7
